$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the 2024-07 row of data (row 51: 油费/过路费/停车费/违章费用/养护费用/顺风车收入/加油退费) ---
$ws.Range("C51").Value = 703
$ws.Range("D51").Value = 183
$ws.Range("E51").Value = 10
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 3365
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0

# Extend the "difference" column formula down into the new row, matching the
# existing per-row formula used in J3:J50 -> (income) - (expenses)
$ws.Range("J51").FormulaR1C1 = "=(RC[-2]+RC[-1])-(RC[-7]+RC[-6]+RC[-5]+RC[-4]+RC[-3])"

# --- Move the selection to where the user left off editing ---
$ws.Range("G52").Select() | Out-Null
